$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.696.01'
$ws.Range('E2').Value = '  -0.72%  '
$ws.Range('D3').Value = '1.895.40'
$ws.Range('E3').Value = '  +1.22%  '
$ws.Range('E4').Value = '  -1.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '312.60'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.44%  '
$ws.Range('E6').Value = '  -1.09%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4898'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.43%  '
$ws.Range('E8').Value = '  -0.76%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07333'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9144'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.83%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.57'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.36%  '
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07684'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.70%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.906.81'
$ws.Range('E13').Value = '  +1.47%  '
$ws.Range('E14').Value = '  -0.47%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.619'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.16%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '91.10'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.25%  '
$ws.Range('E17').Value = '  -1.18%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008780'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.08%  '
$ws.Range('E19').Value = '  -1.03%  '
$ws.Range('D20').Value = '27.933.97'
$ws.Range('E20').Value = '  -0.33%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.50'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.57%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.126'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').Value = '2.141.80'
$ws.Range('E23').Value = '  +1.32%  '
$ws.Range('E24').Value = '  -0.93%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.905'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.18%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '153.68'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.31%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.159'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.52%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '115.74'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.29%  '
$ws.Range('E30').Value = '  -1.77%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08911'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.11%  '
$ws.Range('E32').Value = '  -4.22%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.224'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.11%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7652'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.45%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.642'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.56%  '
$ws.Range('E36').Value = '  -0.85%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.536'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -6.44%  '
$ws.Range('E38').Value = '  -3.45%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05280'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.84%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5486'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.980'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.56%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.918'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.93%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.522'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.86%  '
$ws.Range('E44').Value = '  -1.30%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '111.20'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +5.57%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.61'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.27%  '
$ws.Range('E47').Value = '  -2.02%  '
$ws.Range('E48').Value = '  -1.19%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.633'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.46%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '67.61'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.52%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06051'
$ws.Range('D51').Style = 'Normal'
